# Update "想去人数" (want-to-go count) figures to the latest scraped values.
$wb = $excel.ActiveWorkbook

# Sheet "展览" (exhibitions)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F5").Value = 2923
$ws1.Range("F7").Value = 397

# Sheet "全部类型" (all types, aggregated)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F5").Value = 2923
$ws4.Range("F9").Value = 397
